# Updated cryptos list with GitHub Actions - refresh price/volume data,
# and swap in the latest two coins that entered/left the tracked list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '62.866.88'
$ws.Range('E2').Value = '  +2.15%  '
# Row 3
$ws.Range('D3').Value = '3.473.37'
$ws.Range('E3').Value = '  +2.53%  '
# Row 4
$ws.Range('E4').Value = '  +0.03%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.21'
$ws.Range('E5').Value = '  +1.02%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.18'
$ws.Range('E6').Value = '  +4.73%  '
# Row 8
$ws.Range('E8').Value = '  +1.37%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.66'
$ws.Range('E9').Value = '  -0.71%  '
# Row 10
$ws.Range('E10').Value = '  +2.28%  '
# Row 11
$ws.Range('E11').Value = '  +4.04%  '
# Row 12
$ws.Range('D12').Value = '4.071.39'
$ws.Range('E12').Value = '  +2.58%  '
# Row 13
$ws.Range('E13').Value = '  +5.18%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.128'
$ws.Range('E14').Value = '  +2.32%  '
# Row 15
$ws.Range('D15').Value = '3.466.03'
$ws.Range('E15').Value = '  +1.82%  '
# Row 16
$ws.Range('E16').Value = '  +1.03%  '
# Row 17
$ws.Range('D17').Value = '62.894.29'
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.37'
$ws.Range('E18').Value = '  +4.04%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.44'
$ws.Range('E19').Value = '  +5.99%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '9.30'
$ws.Range('E20').Value = '  +3.42%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '390.27'
$ws.Range('E21').Value = '  -0.23%  '
# Row 22
$ws.Range('E22').Value = '  +2.29%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '74.97'
$ws.Range('E23').Value = '  -0.38%  '
# Row 24
$ws.Range('E24').Value = '  -0.07%  '
# Row 25
$ws.Range('D25').Value = '3.618.76'
$ws.Range('E25').Value = '  +2.56%  '
# Row 26
$ws.Range('E26').Value = '  +2.49%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.180'
$ws.Range('E27').Value = '  -6.56%  '
# Row 28
$ws.Range('E28').Value = '  +6.05%  '
# Row 29
$ws.Range('E29').Value = '  +0.27%  '
# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.21'
$ws.Range('E30').Value = '  +2.10%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.15'
$ws.Range('E31').Value = '  +0.45%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.41'
$ws.Range('E32').Value = '  +0.79%  '
# Row 33
$ws.Range('E33').Value = '  +0.05%  '
# Row 34
$ws.Range('E34').Value = '  +2.33%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.13'
$ws.Range('E35').Value = '  +3.21%  '
# Row 36
$ws.Range('E36').Value = '  +4.88%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '170.81'
$ws.Range('E37').Value = '  +1.47%  '
# Row 38
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.57'
$ws.Range('E38').Value = '  +7.03%  '
# Row 39
$ws.Range('B39').Value = 'EnergySwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '31.33'
$ws.Range('E39').Value = '  +20.10%  '
# Row 40
$ws.Range('D40').Value = '3.514.12'
$ws.Range('E40').Value = '  +2.67%  '
# Row 41
$ws.Range('E41').Value = '  +0.33%  '
# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.802'
$ws.Range('E42').Value = '  +3.01%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.50'
$ws.Range('E43').Value = '  +1.93%  '
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '42.24'
$ws.Range('E44').Value = '  -0.71%  '
# Row 45
$ws.Range('E45').Value = '  +3.58%  '
# Row 46
$ws.Range('E46').Value = '  +3.20%  '
# Row 47
$ws.Range('D47').Value = '2.603.32'
$ws.Range('E47').Value = '  +6.07%  '
# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.57'
$ws.Range('E48').Value = '  +2.72%  '
# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.26'
$ws.Range('E49').Value = '  +11.86%  '
# Row 50
$ws.Range('E50').Value = '  +2.03%  '
# Row 51
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0269'
$ws.Range('E51').Value = '  +2.29%  '
